$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.502.59"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "1.571.60"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -1.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.991"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.794.78"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.562.62"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "27.473.60"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.99%  "
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "1.457.08"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.540"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.813"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("E44").Value = "  +7.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.978"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").Value = "1.707.18"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.85%  "
